$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) B:E
$ws.Range("B2").Value = 1.1088231157342747
$ws.Range("C2").Value = 0.45355212704528869
$ws.Range("D2").Value = 1.4188735565816963
$ws.Range("E2").Value = 0.41317376586636378

# Update row 3 (STR) B:E
$ws.Range("B3").Value = 1.3703464467278372
$ws.Range("C3").Value = 0.95181260226671338
$ws.Range("D3").Value = 1.4712168892325963
$ws.Range("E3").Value = 0.58953414315802211

# Update the selection to match the new state
$ws.Range("B1:E3").Select()
